# DOMA-3100: add formatter convert to number for some columns
#
# The ticket-count columns (processing / completed / canceled / deferred /
# closed / new_or_reopened) on rows 2 and 3 are template placeholder cells
# consumed by a reporting engine. Previously they were plain text
# placeholders; this change appends the ":formatN()" formatter to each
# placeholder so the exported value is rendered as a number, and switches
# the cells' number format from Text ("@") to a plain integer ("0") so
# Excel itself also treats/renders them as numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (d.tickets[i].*) -------------------------------------------------
$ws.Range("C2").Value = "{d.tickets[i].processing:formatN()}"
$ws.Range("D2").Value = "{d.tickets[i].completed:formatN()}"
$ws.Range("E2").Value = "{d.tickets[i].canceled:formatN()}"
$ws.Range("F2").Value = "{d.tickets[i].deferred:formatN()}"
$ws.Range("G2").Value = "{d.tickets[i].closed:formatN()}"
$ws.Range("H2").Value = "{d.tickets[i].new_or_reopened:formatN()}"

# --- Row 3 (d.tickets[i+1].*) -----------------------------------------------
$ws.Range("C3").Value = "{d.tickets[i+1].processing:formatN()}"
$ws.Range("D3").Value = "{d.tickets[i+1].completed:formatN()}"
$ws.Range("E3").Value = "{d.tickets[i+1].canceled:formatN()}"
$ws.Range("F3").Value = "{d.tickets[i+1].deferred:formatN()}"
$ws.Range("G3").Value = "{d.tickets[i+1].closed:formatN()}"
$ws.Range("H3").Value = "{d.tickets[i+1].new_or_reopened:formatN()}"

# Switch those same cells (C:H on rows 2 & 3) from Text format to a number
# format so the values are recognised/displayed as numbers.
$ws.Range("C2:H3").NumberFormat = "0"
